# Updates cryptos list per upstream refresh (GitHub Actions job)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '52.354.33'
$ws.Range("E2").Value = '  +5.77%  '

# Row 3
$ws.Range("D3").Value = '2.817.22'
$ws.Range("E3").Value = '  +7.12%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '118.15'
$ws.Range("E5").Value = '  +6.17%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '337.39'
$ws.Range("E6").Value = '  +4.13%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.541'
$ws.Range("E7").Value = '  +3.84%  '

# Row 8
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.582'
$ws.Range("E9").Value = '  +7.42%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.63'
$ws.Range("E10").Value = '  +8.45%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0867'
$ws.Range("E11").Value = '  +7.37%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.44'
$ws.Range("E12").Value = '  +2.62%  '

# Row 13
$ws.Range("E13").Value = '  +2.64%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.73'
$ws.Range("E14").Value = '  +5.83%  '

# Row 15
$ws.Range("D15").Value = '3.261.99'
$ws.Range("E15").Value = '  +7.26%  '

# Row 16
$ws.Range("D16").Value = '2.821.47'
$ws.Range("E16").Value = '  +7.53%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.895'
$ws.Range("E17").Value = '  +5.66%  '

# Row 18
$ws.Range("D18").Value = '52.360.89'
$ws.Range("E18").Value = '  +5.96%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.30'
$ws.Range("E19").Value = '  +13.39%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.68'
$ws.Range("E20").Value = '  +5.97%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.98'
$ws.Range("E21").Value = '  +4.94%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0984'
$ws.Range("E22").Value = '  +4.46%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '281.38'
$ws.Range("E23").Value = '  +5.00%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.49'
$ws.Range("E24").Value = '  +2.61%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.76'
$ws.Range("E25").Value = '  +9.23%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.18'
$ws.Range("E26").Value = '  +4.95%  '

# Row 27
$ws.Range("E27").Value = '  +0.02%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.32'
$ws.Range("E28").Value = '  +1.77%  '

# Row 29
$ws.Range("E29").Value = '  -0.62%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.144'
$ws.Range("E30").Value = '  +5.37%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '36.05'
$ws.Range("E31").Value = '  +4.78%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.54'
$ws.Range("E32").Value = '  +2.19%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.69'
$ws.Range("E33").Value = '  +3.98%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0831'
$ws.Range("E34").Value = '  +2.77%  '

# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.16'
$ws.Range("E35").Value = '  +6.67%  '

# Row 36
$ws.Range("B36").Value = 'Celestia'
$ws.Range("C36").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.44'
$ws.Range("E36").Value = '  +3.19%  '

# Row 37
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.35'
$ws.Range("E37").Value = '  +8.50%  '

# Row 38
$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.08%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.05'
$ws.Range("E39").Value = '  +3.29%  '

# Row 40
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.77'
$ws.Range("E40").Value = '  +28.80%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0358'
$ws.Range("E41").Value = '  +11.16%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.77'
$ws.Range("E42").Value = '  +8.04%  '

# Row 43
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '128.32'
$ws.Range("E43").Value = '  +0.15%  '

# Row 44
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.36'
$ws.Range("E44").Value = '  +9.87%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.115'
$ws.Range("E45").Value = '  +3.90%  '

# Row 46
$ws.Range("D46").Value = '2.117.96'
$ws.Range("E46").Value = '  +3.58%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.37'
$ws.Range("E47").Value = '  +6.08%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.24'
$ws.Range("E48").Value = '  +3.93%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.60'
$ws.Range("E49").Value = '  +8.09%  '

# Row 50
$ws.Range("B50").Value = 'SEI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.927'
$ws.Range("E50").Value = '  +24.85%  '

# Row 51
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '60.72'
$ws.Range("E51").Value = '  +4.35%  '
